# Apply the change described by the commit:
#   - ToLocation test data value changes from "mia" to "bost" (cell E2 on Sheet1)
#   - The sheet view's selection moves to D11 (with the sheet scrolled so column A
#     is out of view, topLeftCell="B1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Update the shared-string test data: "mia" -> "bost"
$ws.Range("E2").Value = "bost"

# Restore the current selection/view to match the saved workbook state:
# the user had scrolled right one column and left the cursor on D11.
$ws.Range("D11").Select() | Out-Null
